$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1419
$ws.Range("J17").Value = 1517.7778
$ws.Range("L17").Value = 4553.3334
$ws.Range("N17").Value = -4889.3334
$ws.Range("H19").Value = 5019.9165
$ws.Range("I19").Value = 7071.2856
$ws.Range("J19").Value = 2148
$ws.Range("K19").Value = 7071.2856
$ws.Range("L19").Value = 2148
$ws.Range("M19").Value = -6896.2856
$ws.Range("N19").Value = -2498
$ws.Range("H116").Value = 2000
$ws.Range("J116").Value = 2000
$ws.Range("L116").Value = 2000
$ws.Range("N116").Value = -8884
$ws.Range("H129").Value = 22376.043
$ws.Range("J129").Value = 29852.543
$ws.Range("L129").Value = 89557.629
$ws.Range("N129").Value = -99557.629
$ws.Range("H141").Value = 1688.7018
$ws.Range("I141").Value = 1112.4474
$ws.Range("J141").Value = 2841.2104
$ws.Range("K141").Value = 3337.3422
$ws.Range("L141").Value = 8523.6312
$ws.Range("M141").Value = 1842.6578
$ws.Range("N141").Value = -18883.6312

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3260.5251
$ws.Range("I32").Value = 2633.0234
$ws.Range("K32").Value = 2633.0234
$ws.Range("M32").Value = -2346.0234
$ws.Range("H45").Value = 1144.238
$ws.Range("I45").Value = 1054.3334
$ws.Range("J45").Value = 1369
$ws.Range("K45").Value = 1054.3334
$ws.Range("L45").Value = 1369
$ws.Range("M45").Value = -677.3334
$ws.Range("N45").Value = -2123
$ws.Range("H74").Value = 607.32355
$ws.Range("I74").Value = 623.40625
$ws.Range("J74").Value = 350
$ws.Range("K74").Value = 623.40625
$ws.Range("L74").Value = 350
$ws.Range("M74").Value = 250.59375
$ws.Range("N74").Value = -2098
$ws.Range("H77").Value = 607.32355
$ws.Range("I77").Value = 623.40625
$ws.Range("J77").Value = 350
$ws.Range("K77").Value = 3117.03125
$ws.Range("L77").Value = 1750
$ws.Range("M77").Value = 1250.96875
$ws.Range("N77").Value = -10486
$ws.Range("H102").Value = 1860
$ws.Range("I102").Value = 1766.6666
$ws.Range("K102").Value = 1766.6666
$ws.Range("M102").Value = -144.6666
$ws.Range("H135").Value = 40000
$ws.Range("J135").Value = 40000
$ws.Range("L135").Value = 40000
$ws.Range("N135").Value = -50140

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1101.5883
$ws.Range("I94").Value = 718.7
$ws.Range("J94").Value = 1648.5714
$ws.Range("K94").Value = 718.7
$ws.Range("L94").Value = 1648.5714
$ws.Range("M94").Value = -267.7
$ws.Range("N94").Value = -2550.5714
$ws.Range("H100").Value = 6035.75
$ws.Range("J100").Value = 6035.75
$ws.Range("L100").Value = 6035.75
$ws.Range("N100").Value = -8199.75
$ws.Range("H105").Value = 7660
$ws.Range("I105").Value = 7075
$ws.Range("J105").Value = 10000
$ws.Range("K105").Value = 7075
$ws.Range("L105").Value = 10000
$ws.Range("M105").Value = -5328
$ws.Range("N105").Value = -13494

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 456.13043
$ws.Range("I22").Value = 365.86667
$ws.Range("J22").Value = 625.375
$ws.Range("K22").Value = 365.86667
$ws.Range("L22").Value = 625.375
$ws.Range("M22").Value = -15.86667
$ws.Range("N22").Value = -1325.375
$ws.Range("H58").Value = 744.3953
$ws.Range("I58").Value = 557.6129
$ws.Range("K58").Value = 557.6129
$ws.Range("M58").Value = -354.6129
$ws.Range("H99").Value = 2200
$ws.Range("I99").Value = 2200
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2200
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -702
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 2200
$ws.Range("I126").Value = 2200
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6600
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4130
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 744.3953
$ws.Range("I136").Value = 557.6129
$ws.Range("K136").Value = 1672.8387
$ws.Range("M136").Value = 877.1613000000002

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 16666749
$ws.Range("I4").Value = 16666749
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 50000247
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -50000135
$ws.Range("N4").ClearContents()
$ws.Range("H6").Value = 277.08334
$ws.Range("I6").Value = 250
$ws.Range("K6").Value = 750
$ws.Range("M6").Value = -637
$ws.Range("H87").Value = 1745
$ws.Range("I87").Value = 1745
$ws.Range("K87").Value = 5235
$ws.Range("M87").Value = -3987
$ws.Range("H90").Value = 1745
$ws.Range("I90").Value = 1745
$ws.Range("K90").Value = 15705
$ws.Range("M90").Value = -9465
$ws.Range("H131").Value = 1854415.9
$ws.Range("I131").Value = 5418.636
$ws.Range("J131").Value = 2924888
$ws.Range("K131").Value = 16255.908
$ws.Range("L131").Value = 8774664
$ws.Range("M131").Value = -11215.908
$ws.Range("N131").Value = -8784744

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3000
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 3000
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -6244
$ws.Range("H113").Value = 22728180
$ws.Range("I113").Value = 35714956
$ws.Range("J113").Value = 1320
$ws.Range("K113").Value = 35714956
$ws.Range("L113").Value = 1320
$ws.Range("M113").Value = -35712786
$ws.Range("N113").Value = -5660
$ws.Range("H126").Value = 2669.182
$ws.Range("I126").Value = 3477
$ws.Range("J126").Value = 1699.8
$ws.Range("K126").Value = 10431
$ws.Range("L126").Value = 5099.4
$ws.Range("M126").Value = -7961
$ws.Range("N126").Value = -10039.4

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1822.579
$ws.Range("I100").Value = 1863.8462
$ws.Range("J100").Value = 1733.1666
$ws.Range("K100").Value = 1863.8462
$ws.Range("L100").Value = 1733.1666
$ws.Range("M100").Value = -1322.8462
$ws.Range("N100").Value = -2815.1666
$ws.Range("H104").Value = 28254.143
$ws.Range("J104").Value = 28254.143
$ws.Range("L104").Value = 28254.143
$ws.Range("N104").Value = -35242.143
$ws.Range("H122").Value = 3421.0588
$ws.Range("I122").Value = 4082.25
$ws.Range("J122").Value = 2833.3333
$ws.Range("K122").Value = 12246.75
$ws.Range("L122").Value = 8499.999899999999
$ws.Range("M122").Value = -9796.75
$ws.Range("N122").Value = -13399.9999
$ws.Range("H136").Value = 3050.4119
$ws.Range("I136").Value = 3705.0625
$ws.Range("J136").Value = 1947.8422
$ws.Range("K136").Value = 11115.1875
$ws.Range("L136").Value = 5843.5266
$ws.Range("M136").Value = -8565.1875
$ws.Range("N136").Value = -10943.5266

